$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C78").Value = "30min"
$ws.Range("C78").Style = $ws.Range("B78").Style

$ws.Range("C79").Value = "10min"
$ws.Range("C79").Style = $ws.Range("B78").Style

$ws.Range("C80").Value = "10min"
$ws.Range("C80").Style = $ws.Range("B78").Style

$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("C78:C80").Select()
